# "installed capacity graph added"
# Switch the active sheet from "Nodes" to "Col_Name" and move the
# selection on "Col_Name" to I25, then update the two color codes in
# column B (rows 22/23 - Imports/Exports) to their new values.

$wb = $excel.ActiveWorkbook

$colName = $wb.Worksheets.Item("Col_Name")

# Update the "Imports" / "Exports" color codes on Col_Name.
$colName.Range("B22").Value = "#22EB11"
$colName.Range("B23").Value = "#0E25EC"

# Make Col_Name the active sheet (sets tabSelected on it, clears it on
# the previously-active "Nodes" sheet) and select I25 there.
$colName.Activate()
$colName.Range("I25").Select()
